{"js": "const replacements = [\n  [\"775\u00d72=1550\", \"348\u00d74=1392\"],\n  [\"957\u00d75=4785\", \"553\u00d77=3871\"],\n  [\"338\u00d79=3042\", \"458\u00d73=1374\"],\n  [\"180\u00d73=540\", \"116\u00d74=464\"],\n  [\"276\u00d76=1656\", \"919\u00d76=5514\"],\n  [\"258\u00d79=2322\", \"855\u00d73=2565\"],\n  [\"832\u00d72=1664\", \"337\u00d75=1685\"],\n  [\"182\u00d79=1638\", \"167\u00d76=1002\"],\n  [\"993\u00d73=2979\", \"375\u00d79=3375\"],\n  [\"102\u00d76=612\", \"431\u00d77=3017\"],\n  [\"494\u00d73=1482\", \"303\u00d73=909\"],\n  [\"506\u00d75=2530\", \"562\u00d72=1124\"],\n  [\"828\u00d79=7452\", \"172\u00d73=516\"],\n  [\"418\u00d74=1672\", \"815\u00d72=1630\"],\n  [\"504\u00d79=4536\", \"498\u00d77=3486\"],\n  [\"531\u00d76=3186\", \"110\u00d72=220\"],\n  [\"910\u00d78=7280\", \"405\u00d74=1620\"],\n  [\"971\u00d74=3884\", \"388\u00d75=1940\"],\n  [\"353\u00d72=706\", \"733\u00d78=5864\"],\n  [\"403\u00d75=2015\", \"679\u00d74=2716\"],\n  [\"121\u00d72=242\", \"712\u00d78=5696\"],\n  [\"916\u00d78=7328\", \"316\u00d79=2844\"],\n  [\"814\u00d77=5698\", \"405\u00d74=1620\"],\n  [\"634\u00d79=5706\", \"205\u00d75=1025\"],\n  [\"136\u00d78=1088\", \"524\u00d72=1048\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const searchResults = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  if (searchResults.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of searchResults.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"775\u00d72=1550\", \"348\u00d74=1392\"),\n    @(\"957\u00d75=4785\", \"553\u00d77=3871\"),\n    @(\"338\u00d79=3042\", \"458\u00d73=1374\"),\n    @(\"180\u00d73=540\", \"116\u00d74=464\"),\n    @(\"276\u00d76=1656\", \"919\u00d76=5514\"),\n    @(\"258\u00d79=2322\", \"855\u00d73=2565\"),\n    @(\"832\u00d72=1664\", \"337\u00d75=1685\"),\n    @(\"182\u00d79=1638\", \"167\u00d76=1002\"),\n    @(\"993\u00d73=2979\", \"375\u00d79=3375\"),\n    @(\"102\u00d76=612\", \"431\u00d77=3017\"),\n    @(\"494\u00d73=1482\", \"303\u00d73=909\"),\n    @(\"506\u00d75=2530\", \"562\u00d72=1124\"),\n    @(\"828\u00d79=7452\", \"172\u00d73=516\"),\n    @(\"418\u00d74=1672\", \"815\u00d72=1630\"),\n    @(\"504\u00d79=4536\", \"498\u00d77=3486\"),\n    @(\"531\u00d76=3186\", \"110\u00d72=220\"),\n    @(\"910\u00d78=7280\", \"405\u00d74=1620\"),\n    @(\"971\u00d74=3884\", \"388\u00d75=1940\"),\n    @(\"353\u00d72=706\", \"733\u00d78=5864\"),\n    @(\"403\u00d75=2015\", \"679\u00d74=2716\"),\n    @(\"121\u00d72=242\", \"712\u00d78=5696\"),\n    @(\"916\u00d78=7328\", \"316\u00d79=2844\"),\n    @(\"814\u00d77=5698\", \"405\u00d74=1620\"),\n    @(\"634\u00d79=5706\", \"205\u00d75=1025\"),\n    @(\"136\u00d78=1088\", \"524\u00d72=1048\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$oldText, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2) | Out-Null\n}\n"}
